# Update dSF (column F) values for specific rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = 2
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F21").Value = 0
$ws.Range("F25").Value = -2
$ws.Range("F28").Value = 1
